$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45189 to 45190 for rows 2-221
$ws.Range("C2:C221").Value = 45190

# Row 221 gets ht="15" customHeight="1" - set row height explicitly to 15
$ws.Rows.Item(221).RowHeight = 15

# Append new row 222 with data
$ws.Range("A222").Value = "A 44354-2023"
$ws.Range("B222").Value = 45188
$ws.Range("C222").Value = 45190
$ws.Range("D222").Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Range("E222").Value = "HERRLJUNGA"
$ws.Range("G222").Value = 0.7
$ws.Range("H222").Value = 0
$ws.Range("I222").Value = 0
$ws.Range("J222").Value = 0
$ws.Range("K222").Value = 0
$ws.Range("L222").Value = 0
$ws.Range("M222").Value = 0
$ws.Range("N222").Value = 0
$ws.Range("O222").Value = 0
$ws.Range("P222").Value = 0
$ws.Range("Q222").Value = 0

# Copy styles from row 221 to row 222 (B,C get date format style 1, R gets wrap style 2)
$ws.Range("B221:C221").Copy()
$ws.Range("B222:C222").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("R221").Copy()
$ws.Range("R222").PasteSpecial(-4122)  # xlPasteFormats
